# Marks the first 6 API calls as completed ("si") and fills in the
# details of call #6 (which was a stub copy of call #5), per commit
# "creadas las 6 primeras llamadas".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "modificar archivo de configuracion" (patch, same path as
# row 5) becomes its own distinct call: "crear archivo de configuracion"
# (post, to /configuracion).
$ws.Range("A6").Value = "crear archivo de configuracion"
$ws.Range("B6").Value = "/configuracion"
$ws.Range("D6").Value = "post"

# --- Mark the first 6 calls as finished.
$ws.Range("I2").Value = "si"
$ws.Range("I3").Value = "si"
$ws.Range("I4").Value = "si"
$ws.Range("I5").Value = "si"
$ws.Range("I6").Value = "si"

# I5 picks up the underlined-font style (style index 3 in the original
# workbook) to match the other finished entries' emphasis.
$ws.Range("I5").Font.Underline = $true

# Row 5 grows slightly taller to fit the rewrapped text.
$ws.Rows(5).RowHeight = 91.5

# Restore the view to the top of the sheet with B7 selected.
$ws.Range("B7").Select()
